$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: B17/C17/E17 were stored as text ("54446.0","706.0","0") ---
# Convert them to real numbers (54446, 706, 0), matching the target diff.
$ws.Range("B17").Value = 54446
$ws.Range("C17").Value = 706
$ws.Range("E17").Value = 0

# --- Row 18: brand-new row appended after row 17 ---
# A18 holds the date as plain text "2022-01-21" (same convention as the
# other rows in column A, which are text, not real dates) - force text
# via NumberFormat so Excel doesn't auto-convert the string into a date
# serial, then clear the format again so no extra style sticks around.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "2022-01-21"
$ws.Range("A18").ClearFormats()

# B18, C18, D18, F18, G18 stay blank (no value set).

# E18 holds "-2272.0" as text (not a number) - same trick as above.
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2272.0"
$ws.Range("E18").ClearFormats()
